$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.08%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.78%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.676"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.26%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08339"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.86%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.038"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.50%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.98%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.526"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.53%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.992"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.13%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9306"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.34%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1294"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.17%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1973"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.30%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09446"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.57%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03906"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.32%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1060"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.11%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001295"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.04%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006091"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.03%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.442"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.15%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.324"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.57%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1361"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.36%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2403"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-10.35%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04416"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.46%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001247"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.65%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004363"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.81%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.49%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02831"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.28%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05544"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.60%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007803"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.05%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1439"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.21%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008933"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.21%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002138"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.06%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01177"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.64%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007028"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.46%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.60%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003173"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "4.40%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.60%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.60%"
